# "zone name and device names ratio refactoring"
# Updates the export_to_excel / force_extract flag columns (G/H) on the
# service_tables sheet so the zoning-configuration block (rows 99-117,
# "DATA ANALYSIS 7..11") is turned on, a couple of now-superseded rows
# (93, 95, 98) are turned off, and row 64 (fabric devices identification)
# is turned on. Also moves the saved cursor/selection to H99.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("service_tables")
$ws.Activate()

# Row 64 (DATA ANALYSIS 5 - portshow_aggregated): force_extract 0 -> 1
$ws.Range("H64").Value = 1

# Rows that flip export_to_excel off (1 -> 0)
$ws.Range("G93").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("G98").Value = 0

# Rows 99-107,109,110,112,114,117 (DATA ANALYSIS 7-11 zoning/presentation/
# sensor/fabric-stats/raslog blocks): both export_to_excel and
# force_extract flip 0 -> 1
$rows = @(99,100,101,102,103,104,105,106,107,109,110,112,114,117)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = 1
    $ws.Range("H$r").Value = 1
}

# Move the active selection/cursor to H99, matching the saved view state.
$ws.Range("H99").Select()
